# "Found issue with shared stops"
#
# Adds a "Number of shared stops" column to the Meta sheet, sourced from a
# newly-populated "Shared Stops" sheet (rows 7-18 and extra columns C/D, G/H,
# K/L, M/N added), and updates a few window/view settings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook window geometry (cosmetic, matches the new active tab = Meta)
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Left         = -120
$wb.Windows.Item(1).Top          = -120
$wb.Windows.Item(1).Width        = 29040
$wb.Windows.Item(1).Height       = 15990

# ---------------------------------------------------------------------------
# Shared Stops sheet: new rows of data + new helper columns
# ---------------------------------------------------------------------------
$shared = $wb.Worksheets.Item("Shared Stops")

# New rows 7-16 for columns A/B (plain additional stop coordinates)
$ab = @{
    7  = @(64.858142000000001, -147.850549)
    8  = @(64.847294000000005, -147.81287499999999)
    9  = @(64.843765000000005, -147.81287900000001)
    10 = @(64.837947999999997, -147.812816)
    11 = @(64.835809999999995, -147.81517099999999)
    12 = @(64.840632999999997, -147.72513499999999)
}
foreach ($r in $ab.Keys) {
    $shared.Cells.Item($r, 1).Value() = $ab[$r][0]
    $shared.Cells.Item($r, 2).Value() = $ab[$r][1]
}

# New columns C/D rows 3-5 (another stop-coordinate pair)
$cd = @{
    3 = @(64.840632999999997, -147.72513499999999)
    4 = @(64.841434000000007, -147.71964500000001)
    5 = @(64.857146999999998, -147.69340099999999)
}
foreach ($r in $cd.Keys) {
    $shared.Cells.Item($r, 3).Value() = $cd[$r][0]
    $shared.Cells.Item($r, 4).Value() = $cd[$r][1]
}

# New columns G/H rows 3-12
$gh = @{
    3  = @(64.841434000000007, -147.71964500000001)
    4  = @(64.840712999999994, -147.714844)
    5  = @(64.828755000000001, -147.714787)
    6  = @(64.825524999999999, -147.713221)
    7  = @(64.839010000000002, -147.714099)
    8  = @(64.837226999999999, -147.717251)
    9  = @(64.833945,          -147.71713)
    10 = @(64.831175999999999, -147.715881)
    11 = @(64.822730000000007, -147.71213499999999)
    12 = @(64.838575000000006, -147.71881500000001)
}
foreach ($r in $gh.Keys) {
    $shared.Cells.Item($r, 7).Value()  = $gh[$r][0]
    $shared.Cells.Item($r, 8).Value()  = $gh[$r][1]
}

# New columns K/L rows 3-12 (mirrors G/H)
foreach ($r in $gh.Keys) {
    $shared.Cells.Item($r, 11).Value() = $gh[$r][0]
    $shared.Cells.Item($r, 12).Value() = $gh[$r][1]
}

# New/extended columns M/N rows 4-14
$mn = @{
    4  = @(64.835809999999995, -147.81517099999999)
    5  = @(64.841434000000007, -147.71964500000001)
    6  = @(64.856099,          -147.812825)
    7  = @(64.863343999999998, -147.81904299999999)
    8  = @(64.860802000000007, -147.822926)
    9  = @(64.857146999999998, -147.69340099999999)
    10 = @(64.838026999999997, -147.812386)
    11 = @(64.844121999999999, -147.81242700000001)
    12 = @(64.847216000000003, -147.812431)
    13 = @(64.849359000000007, -147.812431)
    14 = @(64.851911999999999, -147.81247400000001)
}
foreach ($r in $mn.Keys) {
    $shared.Cells.Item($r, 13).Value() = $mn[$r][0]
    $shared.Cells.Item($r, 14).Value() = $mn[$r][1]
}

# Extended columns O/P rows 7-18 (continuation of the existing O/P list)
$op = @{
    7  = @(64.858142000000001, -147.850549)
    8  = @(64.847294000000005, -147.81287499999999)
    9  = @(64.843765000000005, -147.81287900000001)
    10 = @(64.837947999999997, -147.812816)
    11 = @(64.835809999999995, -147.81517099999999)
    12 = @(64.856099,          -147.812825)
    13 = @(64.863343999999998, -147.81904299999999)
    14 = @(64.860802000000007, -147.822926)
    15 = @(64.838026999999997, -147.812386)
    16 = @(64.844121999999999, -147.81242700000001)
    17 = @(64.847216000000003, -147.812431)
    18 = @(64.849359000000007, -147.812431)
}
foreach ($r in $op.Keys) {
    $shared.Cells.Item($r, 15).Value() = $op[$r][0]
    $shared.Cells.Item($r, 16).Value() = $op[$r][1]
}

# View: zoomed to 100%, selection moved, no longer the tab shown on open
$shared.Application.ActiveWindow.Zoom = 100
$shared.Range("M8").Select()

# ---------------------------------------------------------------------------
# All Stops sheet: view settings only (zoom / scroll / selection)
# ---------------------------------------------------------------------------
$allStops = $wb.Worksheets.Item("All Stops")
$allStops.Select()
$allStops.Application.ActiveWindow.Zoom = 100
$allStops.Application.ActiveWindow.ScrollRow = 33
$allStops.Range("M60").Select()

# ---------------------------------------------------------------------------
# Meta sheet: new "Number of shared stops" column
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta")
$meta.Select()

$meta.Cells.Item(1, 4).Value() = "Number of shared stops"

$meta.Cells.Item(2, 4).Formula() = "=COUNT('Shared Stops'!A`$3:A`$235)"
$meta.Cells.Item(3, 4).Formula() = "=COUNT('Shared Stops'!C`$3:C`$235)"
$meta.Cells.Item(4, 4).Formula() = "=COUNT('Shared Stops'!E`$3:E`$235)"
$meta.Cells.Item(5, 4).Formula() = "=COUNT('Shared Stops'!G`$3:G`$235)"
$meta.Cells.Item(6, 4).Formula() = "=COUNT('Shared Stops'!I`$3:I`$235)"
$meta.Cells.Item(7, 4).Formula() = "=COUNT('Shared Stops'!K`$3:K`$235)"
$meta.Cells.Item(8, 4).Formula() = "=COUNT('Shared Stops'!M`$3:M`$235)"
$meta.Cells.Item(9, 4).Formula() = "=COUNT('Shared Stops'!O`$3:O`$235)"

$meta.Range("D2:D9").HorizontalAlignment = -4108

$meta.Columns.Item(4).ColumnWidth = 22.5703125

# B2 formula gains an absolute column reference
$meta.Cells.Item(2, 2).Formula() = "=COUNT('Raw stops'!`$A3:`$A235)"

$meta.Range("D7").Select()
